$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Pediatric Vaccine ")
$ws.Range("A2").Value = "DTaP "
$ws.Range("A3").Value = "DTaP "
$ws.Range("A4").Value = "DTaP "
$ws.Range("A5").Value = "DTaP-IPV "
$ws.Range("A6").Value = "DTaP-IPV "
$ws.Range("A7").Value = "DTaP-Hep B-IPV "
$ws.Range("A8").Value = "DTaP-IP-HI "
$ws.Range("A9").Value = "e-IPV "
$ws.Range("A10").Value = "Hepatitis A Pediatric "
$ws.Range("A11").Value = "Hepatitis A Pediatric "
$ws.Range("A12").Value = "Hepatitis A Pediatric "
$ws.Range("A13").Value = "Hepatitis A Pediatric "
$ws.Range("A14").Value = "Hepatitis A-Hepatitis B 18 only "
$ws.Range("A15").Value = "Hepatitis A-Hepatitis B 18 only "
$ws.Range("A16").Value = "Hepatitis B  Pediatric/Adolescent"
$ws.Range("A17").Value = "Hepatitis B  Pediatric/Adolescent"
$ws.Range("A18").Value = "Hepatitis B  Pediatric/Adolescent"
$ws.Range("B18").Value = "Recombivax HB"
$ws.Range("A19").Value = "Hepatitis B  Pediatric/Adolescent"
$ws.Range("B19").Value = "Recombivax HB"
$ws.Range("A20").Value = "Hib "
$ws.Range("A21").Value = "Hib "
$ws.Range("A22").Value = "HIBMENCY "
$ws.Range("A23").Value = "HPV - Quadrivalent Human Papillomavirus Types 6, 11, 16 and 18 Recombinant "
$ws.Range("A24").Value = "HPV - Human Papillomavirus 9-valent "
$ws.Range("A25").Value = "HPV -Bivalent Human Papillomavirus Types 16 and 18 "
$ws.Range("A26").Value = "MENB - Meningococcal Group B "
$ws.Range("A27").Value = "MENB - Meningococcal Group B "
$ws.Range("A28").Value = "MENB - Meningococcal Group B "
$ws.Range("A29").Value = "Meningococcal Conjugate (Groups A, C, Y and W-135) "
$ws.Range("A30").Value = "Meningococcal Conjugate (Groups A, C, Y and W-135) "
$ws.Range("A31").Value = "Measles, Mumps and Rubella (MMR) "
$ws.Range("A32").Value = "MMR/Varicella "
$ws.Range("A33").Value = "Pneumococcal 13-valent  (Pediatric)"
$ws.Range("A35").Value = "Rotavirus, Live, Oral, Pentavalent "
$ws.Range("A36").Value = "Rotavirus, Live, Oral, Pentavalent "
$ws.Range("A37").Value = "Rotavirus, Live, Oral, Oral "
$ws.Range("A38").Value = "Tetanus  Diphtheria Toxoids "
$ws.Range("A39").Value = "Tetanus  Diphtheria Toxoids "
$ws.Range("A40").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws.Range("A41").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws.Range("A42").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws.Range("A43").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws.Range("A44").Value = "Varicella "

$ws = $wb.Worksheets.Item("Adult Vaccine ")
$ws.Range("A2").Value = "Hepatitis A-Adult "
$ws.Range("A3").Value = "Hepatitis A-Adult "
$ws.Range("A4").Value = "Hepatitis A Adult "
$ws.Range("A5").Value = "Hepatitis A Adult "
$ws.Range("A6").Value = "Hepatitis A-Hepatitis B Adult "
$ws.Range("A7").Value = "Hepatitis A-Hepatitis B Adult "
$ws.Range("A8").Value = "Hepatitis B-Adult "
$ws.Range("A9").Value = "Hepatitis B-Adult "
$ws.Range("A10").Value = "Hepatitis B-Adult "
$ws.Range("A11").Value = "Hepatitis B-Adult "
$ws.Range("A12").Value = "HPV -Quadrivalent Human Papillomavirus Types 6, 11, 16 and 18 Recombinant Adult "
$ws.Range("A13").Value = "HPV-Human Papillomavirus 9 Valent "
$ws.Range("A14").Value = "HPV-Human Papillomavirus Bivalent Types 16 and 18 "
$ws.Range("A15").Value = "Measles, Mumps,  Rubella-Adult "
$ws.Range("A16").Value = "Meningococcal Conjugate "
$ws.Range("A17").Value = "Meningococcal Conjugate "
$ws.Range("A18").Value = "MENB - Meningococcal Group B "
$ws.Range("A19").Value = "MENB - Meningococcal Group B "
$ws.Range("A20").Value = "MENB - Meningococcal Group B "
$ws.Range("A21").Value = "Pneumococcal 13-valent  (Adult)"
$ws.Range("A24").Value = "Tetanus and Diphtheria Toxoids"
$ws.Range("A25").Value = "Tetanus and Diphtheria Toxoids"
$ws.Range("A26").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws.Range("A27").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws.Range("A28").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws.Range("A29").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws.Range("A30").Value = "Varicella-Adult "

$ws = $wb.Worksheets.Item("Pediatric Influenza Vaccine ")
$ws.Range("A2").Value = "Influenza  (Age 6 months and older)"
$ws.Range("A3").Value = "Influenza  (Age 6 months and older)"
$ws.Range("B3").Value = "Fluzone Quadrivalent"
$ws.Range("A4").Value = "Influenza  (Age 6-35 months)"
$ws.Range("B4").Value = "Fluzone Quadrivalent Pediatric dose No Preservative"
$ws.Range("A5").Value = "Influenza  (Age 36 months and older)"
$ws.Range("B5").Value = "Fluzone Quadrivalent No-Preservative"
$ws.Range("A6").Value = "Influenza  (Age 36 months and older)"
$ws.Range("B6").Value = "Fluzone Quadrivalent No-Preservative"
$ws.Range("A7").Value = "Influenza  (Age 36 months and older)"
$ws.Range("B7").Value = "Fluarix Quadrivalent Preservative Free"
$ws.Range("D7").Value = "10 pack- 1 dose TipLok syringe"
$ws.Range("A8").Value = "Influenza  (Age 36 months and older)"
$ws.Range("B8").Value = "FluLaval Quadrivalent"
$ws.Range("A9").Value = "Influenza  (Age 4 years and older)"
$ws.Range("A10").Value = "Influenza  (Age 4 years and older)"
$ws.Range("B10").Value = "Fluvirin Preservative Free"
$ws.Range("A11").Value = "Influenza  Live, Intranasal (Age 2-49 years)"
$ws.Range("B11").Value = "FluMist Quadrivalent No Preservative"
$ws.Range("A12").Value = "Influenza  (Age 9 years and older)"
$ws.Range("B12").Value = "Afluria No Preservative"
$ws.Range("D12").Value = "10 pack-1 dose syringe"
$ws.Range("A13").Value = "Influenza  (Age 9 years and older)"

$ws = $wb.Worksheets.Item("Adult Influenza Vaccine ")
$ws.Range("A2").Value = "Influenza  (Age 6 months and older)"
$ws.Range("A3").Value = "Influenza  (Age 6 months and older)"
$ws.Range("B3").Value = "Fluzone Quadrivalent"
$ws.Range("A4").Value = "Influenza  (Age 36 months and older)"
$ws.Range("B4").Value = "Fluzone Quadrivalent No Preservative"
$ws.Range("A5").Value = "Influenza  (Age 36 months and older)"
$ws.Range("B5").Value = "Fluzone Quadrivalent No Preservative"
$ws.Range("A6").Value = "Influenza  (Age 18 years and older)"
$ws.Range("B6").Value = "Flucelvax Preservative Free Antibiotic free"
$ws.Range("A7").Value = "Influenza  (Age 4 years and older)"
$ws.Range("A8").Value = "Influenza  (Age 4 years and older)"
$ws.Range("B8").Value = "Fluvirin Preservative Free"
$ws.Range("A9").Value = "Influenza  (Age 36 months and older)"
$ws.Range("B9").Value = "Fluarix Quadrivalent Preservative Free"
$ws.Range("D9").Value = "10 pack- 1 dose TipLok syringe"
$ws.Range("A10").Value = "Influenza  (Age 36 months and older)"
$ws.Range("A11").Value = "Influenza  (Age 9 years and older)"
$ws.Range("B11").Value = "Afluria No Preservative"
$ws.Range("D11").Value = "10 pack-1 dose syringe"
$ws.Range("A12").Value = "Influenza  (Age 9 years and older)"
$ws.Range("A13").Value = "Influenza  (Age 18 years and older)"
